$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the source diff: CellRef / NewValue / ForceText.
# ForceText=$true is used for column D ("Price") because many of its values look
# like plain numbers (e.g. "1.005", "302.20") but must be stored as literal text
# (matching the original inlineStr cells) rather than being auto-converted to a
# numeric Value by Excel, which would lose formatting such as trailing zeros.
$updates = @(
    @{ Cell = 'D2'; Value = '23.273.21'; ForceText = $true },
    @{ Cell = 'E2'; Value = '  -2.57%  '; ForceText = $false },
    @{ Cell = 'D3'; Value = '1.598.28'; ForceText = $true },
    @{ Cell = 'E3'; Value = '  -3.47%  '; ForceText = $false },
    @{ Cell = 'D4'; Value = '1.005'; ForceText = $true },
    @{ Cell = 'E4'; Value = '  +0.21%  '; ForceText = $false },
    @{ Cell = 'D5'; Value = '1.005'; ForceText = $true },
    @{ Cell = 'E5'; Value = '  +0.23%  '; ForceText = $false },
    @{ Cell = 'D6'; Value = '302.20'; ForceText = $true },
    @{ Cell = 'E6'; Value = '  -2.33%  '; ForceText = $false },
    @{ Cell = 'D7'; Value = '0.3761'; ForceText = $true },
    @{ Cell = 'E7'; Value = '  -3.36%  '; ForceText = $false },
    @{ Cell = 'D8'; Value = '0.3659'; ForceText = $true },
    @{ Cell = 'E8'; Value = '  -4.79%  '; ForceText = $false },
    @{ Cell = 'D9'; Value = '49.74'; ForceText = $true },
    @{ Cell = 'E9'; Value = '  -2.40%  '; ForceText = $false },
    @{ Cell = 'D10'; Value = '1.005'; ForceText = $true },
    @{ Cell = 'E10'; Value = '  +0.28%  '; ForceText = $false },
    @{ Cell = 'D11'; Value = '1.276'; ForceText = $true },
    @{ Cell = 'E11'; Value = '  -5.77%  '; ForceText = $false },
    @{ Cell = 'D12'; Value = '0.08113'; ForceText = $true },
    @{ Cell = 'E12'; Value = '  -4.14%  '; ForceText = $false },
    @{ Cell = 'D13'; Value = '22.88'; ForceText = $true },
    @{ Cell = 'E13'; Value = '  -4.23%  '; ForceText = $false },
    @{ Cell = 'D14'; Value = '6.649'; ForceText = $true },
    @{ Cell = 'E14'; Value = '  -7.28%  '; ForceText = $false },
    @{ Cell = 'D15'; Value = '7.569'; ForceText = $true },
    @{ Cell = 'E15'; Value = '  -4.21%  '; ForceText = $false },
    @{ Cell = 'D16'; Value = '0.00001267'; ForceText = $true },
    @{ Cell = 'E16'; Value = '  -2.83%  '; ForceText = $false },
    @{ Cell = 'D17'; Value = '1.598.43'; ForceText = $true },
    @{ Cell = 'E17'; Value = '  -3.40%  '; ForceText = $false },
    @{ Cell = 'D18'; Value = '91.42'; ForceText = $true },
    @{ Cell = 'E18'; Value = '  -3.56%  '; ForceText = $false },
    @{ Cell = 'D19'; Value = '0.06825'; ForceText = $true },
    @{ Cell = 'E19'; Value = '  -2.51%  '; ForceText = $false },
    @{ Cell = 'D20'; Value = '18.51'; ForceText = $true },
    @{ Cell = 'E20'; Value = '  -6.59%  '; ForceText = $false },
    @{ Cell = 'D21'; Value = '6.620'; ForceText = $true },
    @{ Cell = 'E21'; Value = '  -4.24%  '; ForceText = $false },
    @{ Cell = 'E22'; Value = '  +0.18%  '; ForceText = $false },
    @{ Cell = 'D23'; Value = '13.15'; ForceText = $true },
    @{ Cell = 'E23'; Value = '  -3.22%  '; ForceText = $false },
    @{ Cell = 'D24'; Value = '23.264.19'; ForceText = $true },
    @{ Cell = 'E24'; Value = '  -2.63%  '; ForceText = $false },
    @{ Cell = 'D25'; Value = '2.384'; ForceText = $true },
    @{ Cell = 'E25'; Value = '  -4.44%  '; ForceText = $false },
    @{ Cell = 'D26'; Value = '2.965'; ForceText = $true },
    @{ Cell = 'E26'; Value = '  -2.51%  '; ForceText = $false },
    @{ Cell = 'D27'; Value = '21.16'; ForceText = $true },
    @{ Cell = 'E27'; Value = '  -4.22%  '; ForceText = $false },
    @{ Cell = 'D28'; Value = '150.56'; ForceText = $true },
    @{ Cell = 'E28'; Value = '  -1.59%  '; ForceText = $false },
    @{ Cell = 'D29'; Value = '5.321'; ForceText = $true },
    @{ Cell = 'E29'; Value = '  -1.46%  '; ForceText = $false },
    @{ Cell = 'D30'; Value = '132.49'; ForceText = $true },
    @{ Cell = 'E30'; Value = '  -5.04%  '; ForceText = $false },
    @{ Cell = 'D31'; Value = '2.477'; ForceText = $true },
    @{ Cell = 'E31'; Value = '  -0.56%  '; ForceText = $false },
    @{ Cell = 'D32'; Value = '7.152'; ForceText = $true },
    @{ Cell = 'E32'; Value = '  -8.62%  '; ForceText = $false },
    @{ Cell = 'D33'; Value = '1.772.53'; ForceText = $true },
    @{ Cell = 'E33'; Value = '  -3.48%  '; ForceText = $false },
    @{ Cell = 'D34'; Value = '0.9630'; ForceText = $true },
    @{ Cell = 'E34'; Value = '  -6.98%  '; ForceText = $false },
    @{ Cell = 'D35'; Value = '0.07724'; ForceText = $true },
    @{ Cell = 'E35'; Value = '  -4.27%  '; ForceText = $false },
    @{ Cell = 'D36'; Value = '0.02785'; ForceText = $true },
    @{ Cell = 'E36'; Value = '  -6.09%  '; ForceText = $false },
    @{ Cell = 'D37'; Value = '6.289'; ForceText = $true },
    @{ Cell = 'E37'; Value = '  -5.75%  '; ForceText = $false },
    @{ Cell = 'D38'; Value = '10.27'; ForceText = $true },
    @{ Cell = 'E38'; Value = '  -6.95%  '; ForceText = $false },
    @{ Cell = 'D39'; Value = '0.2550'; ForceText = $true },
    @{ Cell = 'E39'; Value = '  -5.13%  '; ForceText = $false },
    @{ Cell = 'D40'; Value = '0.08870'; ForceText = $true },
    @{ Cell = 'E40'; Value = '  -2.79%  '; ForceText = $false },
    @{ Cell = 'D41'; Value = '1.389'; ForceText = $true },
    @{ Cell = 'E41'; Value = '  -2.07%  '; ForceText = $false },
    @{ Cell = 'B42'; Value = 'Aptos'; ForceText = $false },
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; ForceText = $false },
    @{ Cell = 'D42'; Value = '12.86'; ForceText = $true },
    @{ Cell = 'E42'; Value = '  -4.85%  '; ForceText = $false },
    @{ Cell = 'B43'; Value = 'TheSandbox'; ForceText = $false },
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; ForceText = $false },
    @{ Cell = 'D43'; Value = '0.7175'; ForceText = $true },
    @{ Cell = 'E43'; Value = '  -4.71%  '; ForceText = $false },
    @{ Cell = 'D44'; Value = '16.11'; ForceText = $true },
    @{ Cell = 'E44'; Value = '  -2.19%  '; ForceText = $false },
    @{ Cell = 'D45'; Value = '0.6628'; ForceText = $true },
    @{ Cell = 'E45'; Value = '  -4.75%  '; ForceText = $false },
    @{ Cell = 'D46'; Value = '2.315'; ForceText = $true },
    @{ Cell = 'E46'; Value = '  -6.12%  '; ForceText = $false },
    @{ Cell = 'D47'; Value = '1.004'; ForceText = $true },
    @{ Cell = 'E47'; Value = '  +0.17%  '; ForceText = $false },
    @{ Cell = 'D48'; Value = '3.976'; ForceText = $true },
    @{ Cell = 'E48'; Value = '  -2.48%  '; ForceText = $false },
    @{ Cell = 'D49'; Value = '0.08013'; ForceText = $true },
    @{ Cell = 'E49'; Value = '  -3.06%  '; ForceText = $false },
    @{ Cell = 'D50'; Value = '131.89'; ForceText = $true },
    @{ Cell = 'E50'; Value = '  -2.04%  '; ForceText = $false },
    @{ Cell = 'D51'; Value = '1.177'; ForceText = $true },
    @{ Cell = 'E51'; Value = '  -4.86%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}

